$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "UDES"
$ws.Range("B3").Value = "asdasd"
$ws.Range("C3").Value = "otros"
$ws.Range("D3").Value = "24 horas"
$ws.Range("E3").Value = 4

# "123123123" and "123123" look like plain numbers, so a direct .Value
# assignment would store them as numeric cells. Enter them as text formulas
# first, then convert the formula results to static values via copy / paste
# special (values only) so the cells end up as plain shared-string text
# cells without requiring any new (text) number-format style.
$ws.Range("F3").Formula = '="123123123"'
$ws.Range("G3").Formula = '="123123"'
$ws.Range("F3:G3").Copy()
$ws.Range("F3:G3").PasteSpecial(-4163)

$ws.Range("H3").Value = "ghola soy javier"
